# Apply scheduled market-data refresh to the Moogle Profits workbook.
# For each changed cell: update its value; for removed cells, clear them;
# for newly-populated cells, set their value (creating the cell).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1303.6364  # was 1391
$ws.Range("I19").Value = 1607.3334  # was 3001
$ws.Range("K19").Value = 1607.3334  # was 3001
$ws.Range("M19").Value = -1432.3334  # was -2826
$ws.Range("H31").Value = 7097.3335  # was 1757.7142
$ws.Range("I31").Value = 303  # was 217.5
$ws.Range("J31").Value = 10494.5  # was 10999
$ws.Range("K31").Value = 909  # was 652.5
$ws.Range("L31").Value = 31483.5  # was 32997
$ws.Range("M31").Value = -679  # was -422.5
$ws.Range("N31").Value = -31943.5  # was -33457
$ws.Range("H62").Value = 5125  # was 3573
$ws.Range("I62").Value = 5125  # was 3834.1667
$ws.Range("J62").Value = 0  # was 2006
$ws.Range("K62").Value = 5125  # was 3834.1667
$ws.Range("L62").Value = 0  # was 2006
$ws.Range("M62").Value = -4501  # was -3210.1667
$ws.Range("N62").ClearContents()  # was -3254
$ws.Range("H65").Value = 5125  # was 3573
$ws.Range("I65").Value = 5125  # was 3834.1667
$ws.Range("J65").Value = 0  # was 2006
$ws.Range("K65").Value = 25625  # was 19170.8335
$ws.Range("L65").Value = 0  # was 10030
$ws.Range("M65").Value = -22505  # was -16050.8335
$ws.Range("N65").ClearContents()  # was -16270
$ws.Range("H116").Value = 6225  # was 4544.875
$ws.Range("I116").Value = 9900  # was 4227.25
$ws.Range("J116").Value = 5000  # was 4862.5
$ws.Range("K116").Value = 9900  # was 4227.25
$ws.Range("L116").Value = 5000  # was 4862.5
$ws.Range("M116").Value = -6458  # was -785.25
$ws.Range("N116").Value = -11884  # was -11746.5
$ws.Range("H141").Value = 6466.6665  # was 5222.154
$ws.Range("I141").Value = 4111.1113  # was 3788.8
$ws.Range("K141").Value = 12333.3339  # was 11366.4
$ws.Range("M141").Value = -7153.333899999998  # was -6186.400000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2479.4  # was 2439.0667
$ws.Range("I45").Value = 1432.7778  # was 1379
$ws.Range("J45").Value = 4049.3333  # was 4559.2
$ws.Range("K45").Value = 1432.7778  # was 1379
$ws.Range("L45").Value = 4049.3333  # was 4559.2
$ws.Range("M45").Value = -1055.7778  # was -1002
$ws.Range("N45").Value = -4803.3333  # was -5313.2
$ws.Range("H61").Value = 4573.3335  # was 4742.857
$ws.Range("J61").Value = 7702.857  # was 12666.667
$ws.Range("L61").Value = 7702.857  # was 12666.667
$ws.Range("N61").Value = -8126.857  # was -13090.667
$ws.Range("H88").Value = 3734  # was 3839.3333
$ws.Range("J88").Value = 3715  # was 3850.4285
$ws.Range("L88").Value = 3715  # was 3850.4285
$ws.Range("N88").Value = -4527  # was -4662.4285
$ws.Range("H91").Value = 3734  # was 3839.3333
$ws.Range("J91").Value = 3715  # was 3850.4285
$ws.Range("L91").Value = 3715  # was 3850.4285
$ws.Range("N91").Value = -6523  # was -6658.4285
$ws.Range("H136").Value = 4573.3335  # was 4742.857
$ws.Range("J136").Value = 7702.857  # was 12666.667
$ws.Range("L136").Value = 23108.571  # was 38000.001
$ws.Range("N136").Value = -28208.571  # was -43100.001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 19500  # was 25000
$ws.Range("I14").Value = 19500  # was 25000
$ws.Range("K14").Value = 19500  # was 25000
$ws.Range("M14").Value = -19328  # was -24828
$ws.Range("H20").Value = 2156.6  # was 2106.8
$ws.Range("I20").Value = 927.8333  # was 1009.5714
$ws.Range("J20").Value = 3999.75  # was 4667
$ws.Range("K20").Value = 927.8333  # was 1009.5714
$ws.Range("L20").Value = 3999.75  # was 4667
$ws.Range("M20").Value = -680.8333  # was -762.5714
$ws.Range("N20").Value = -4493.75  # was -5161
$ws.Range("H58").Value = 23646.334  # was 23609.25
$ws.Range("J58").Value = 23646.334  # was 23609.25
$ws.Range("L58").Value = 23646.334  # was 23609.25
$ws.Range("N58").Value = -24234.334  # was -24197.25
$ws.Range("H134").Value = 4009.8408  # was 3949.3555
$ws.Range("I134").Value = 2042.8529  # was 2045.7354
$ws.Range("J134").Value = 10697.6  # was 9833.272000000001
$ws.Range("K134").Value = 6128.5587  # was 6137.206200000001
$ws.Range("L134").Value = 32092.8  # was 29499.816
$ws.Range("M134").Value = -3593.5587  # was -3602.206200000001
$ws.Range("N134").Value = -37162.8  # was -34569.81600000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 72999.5  # was 76633
$ws.Range("I52").Value = 60000  # was 71950
$ws.Range("K52").Value = 60000  # was 71950
$ws.Range("M52").Value = -59706  # was -71656
$ws.Range("H59").Value = 94142.14  # was 98165.836
$ws.Range("J59").Value = 103165.836  # was 109799
$ws.Range("L59").Value = 103165.836  # was 109799
$ws.Range("N59").Value = -105455.836  # was -112089

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 417.875  # was 399.17648
$ws.Range("I2").Value = 627  # was 561.125
$ws.Range("K2").Value = 3762  # was 3366.75
$ws.Range("M2").Value = -3649  # was -3253.75
$ws.Range("H12").Value = 114  # was 104.6
$ws.Range("I12").Value = 1  # was 10.5
$ws.Range("K12").Value = 3  # was 31.5
$ws.Range("M12").Value = 170  # was 141.5
$ws.Range("H38").Value = 69.5  # was 68.22221999999999
$ws.Range("I38").Value = 63.714287  # was 68.22221999999999
$ws.Range("J38").Value = 89.75  # was 0
$ws.Range("K38").Value = 191.142861  # was 204.66666
$ws.Range("L38").Value = 269.25  # was 0
$ws.Range("M38").Value = 155.857139  # was 142.33334
$ws.Range("N38").Value = -963.25  # newly populated
$ws.Range("H54").Value = 1999.6666  # was 5336
$ws.Range("I54").Value = 2004  # was 3004
$ws.Range("J54").Value = 1997.5  # was 10000
$ws.Range("K54").Value = 6012  # was 9012
$ws.Range("L54").Value = 5992.5  # was 30000
$ws.Range("M54").Value = -5453  # was -8453
$ws.Range("N54").Value = -7110.5  # was -31118
$ws.Range("H57").Value = 14999.857  # was 12444.223
$ws.Range("I57").Value = 7000  # was 5999.6665
$ws.Range("J57").Value = 18199.8  # was 15666.5
$ws.Range("K57").Value = 21000  # was 17998.9995
$ws.Range("L57").Value = 54599.39999999999  # was 46999.5
$ws.Range("M57").Value = -20441  # was -17439.9995
$ws.Range("N57").Value = -55717.39999999999  # was -48117.5
$ws.Range("H58").Value = 12547  # was 12647.25
$ws.Range("I58").Value = 12547  # was 12647.25
$ws.Range("K58").Value = 37641  # was 37941.75
$ws.Range("M58").Value = -37513  # was -37813.75
$ws.Range("H86").Value = 942.63635  # was 1007
$ws.Range("J86").Value = 1187.8  # was 1410
$ws.Range("L86").Value = 3563.4  # was 4230
$ws.Range("N86").Value = -5935.4  # was -6602
$ws.Range("H89").Value = 942.63635  # was 1007
$ws.Range("J89").Value = 1187.8  # was 1410
$ws.Range("L89").Value = 10690.2  # was 12690
$ws.Range("N89").Value = -22546.2  # was -24546
$ws.Range("H98").Value = 397.125  # was 440.25
$ws.Range("I98").Value = 378.33334  # was 410.33334
$ws.Range("J98").Value = 408.4  # was 450.22223
$ws.Range("K98").Value = 1135.00002  # was 1231.00002
$ws.Range("L98").Value = 1225.2  # was 1350.66669
$ws.Range("M98").Value = 362.9999800000001  # was 266.9999800000001
$ws.Range("N98").Value = -4221.2  # was -4346.66669
$ws.Range("H116").Value = 7369.25  # was 8119.25
$ws.Range("J116").Value = 6493.5  # was 7993.5
$ws.Range("L116").Value = 19480.5  # was 23980.5
$ws.Range("N116").Value = -26364.5  # was -30864.5
$ws.Range("H141").Value = 6270.8125  # was 6270.875
$ws.Range("I141").Value = 5015.273  # was 5015.364
$ws.Range("K141").Value = 15045.819  # was 15046.092
$ws.Range("M141").Value = -9865.819  # was -9866.091999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7445.625  # was 7595.4
$ws.Range("J70").Value = 8366.166999999999  # was 8999.6
$ws.Range("L70").Value = 8366.166999999999  # was 8999.6
$ws.Range("N70").Value = -8906.166999999999  # was -9539.6
$ws.Range("H73").Value = 7445.625  # was 7595.4
$ws.Range("J73").Value = 8366.166999999999  # was 8999.6
$ws.Range("L73").Value = 8366.166999999999  # was 8999.6
$ws.Range("N73").Value = -10238.167  # was -10871.6
$ws.Range("H92").Value = 155128.25  # was 36302.8
$ws.Range("J92").Value = 155128.25  # was 36302.8
$ws.Range("L92").Value = 155128.25  # was 36302.8
$ws.Range("N92").Value = -158872.25  # was -40046.8
$ws.Range("H132").Value = 3231.102  # was 3345.1702
$ws.Range("I132").Value = 2645.3865  # was 2709.186
$ws.Range("J132").Value = 8385.4  # was 10182
$ws.Range("K132").Value = 7936.1595  # was 8127.558000000001
$ws.Range("L132").Value = 25156.2  # was 30546
$ws.Range("M132").Value = -5406.1595  # was -5597.558000000001
$ws.Range("N132").Value = -30216.2  # was -35606

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2725.3635  # was 3053.3333
$ws.Range("I22").Value = 2717.375  # was 2962.8572
$ws.Range("J22").Value = 2746.6667  # was 3370
$ws.Range("K22").Value = 2717.375  # was 2962.8572
$ws.Range("L22").Value = 2746.6667  # was 3370
$ws.Range("M22").Value = -2422.375  # was -2667.8572
$ws.Range("N22").Value = -3336.6667  # was -3960
$ws.Range("H27").Value = 2725.3635  # was 3053.3333
$ws.Range("I27").Value = 2717.375  # was 2962.8572
$ws.Range("J27").Value = 2746.6667  # was 3370
$ws.Range("K27").Value = 2717.375  # was 2962.8572
$ws.Range("L27").Value = 2746.6667  # was 3370
$ws.Range("M27").Value = -2610.375  # was -2855.8572
$ws.Range("N27").Value = -2960.6667  # was -3584

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12070  # was 11887.556
$ws.Range("I62").Value = 7750  # was 8199.799999999999
$ws.Range("J62").Value = 17830  # was 16497.25
$ws.Range("K62").Value = 7750  # was 8199.799999999999
$ws.Range("L62").Value = 17830  # was 16497.25
$ws.Range("M62").Value = -7126  # was -7575.799999999999
$ws.Range("N62").Value = -19078  # was -17745.25
$ws.Range("H65").Value = 12070  # was 11887.556
$ws.Range("I65").Value = 7750  # was 8199.799999999999
$ws.Range("J65").Value = 17830  # was 16497.25
$ws.Range("K65").Value = 38750  # was 40999
$ws.Range("L65").Value = 89150  # was 82486.25
$ws.Range("M65").Value = -35630  # was -37879
$ws.Range("N65").Value = -95390  # was -88726.25
$ws.Range("H81").Value = 4434  # was 4328.3335
$ws.Range("I81").Value = 3339  # was 3308.182
$ws.Range("J81").Value = 5650.6665  # was 5450.5
$ws.Range("K81").Value = 6678  # was 6616.364
$ws.Range("L81").Value = 11301.333  # was 10901
$ws.Range("M81").Value = -5617  # was -5555.364
$ws.Range("N81").Value = -13423.333  # was -13023
$ws.Range("H84").Value = 4434  # was 4328.3335
$ws.Range("I84").Value = 3339  # was 3308.182
$ws.Range("J84").Value = 5650.6665  # was 5450.5
$ws.Range("K84").Value = 33390  # was 33081.82
$ws.Range("L84").Value = 56506.665  # was 54505
$ws.Range("M84").Value = -28086  # was -27777.82
$ws.Range("N84").Value = -67114.66500000001  # was -65113
$ws.Range("H100").Value = 1003  # was 990.36
$ws.Range("I100").Value = 847.7857  # was 810.9286
$ws.Range("J100").Value = 1200.5454  # was 1218.7273
$ws.Range("K100").Value = 1695.5714  # was 1621.8572
$ws.Range("L100").Value = 2401.0908  # was 2437.4546
$ws.Range("M100").Value = -1154.5714  # was -1080.8572
$ws.Range("N100").Value = -3483.0908  # was -3519.4546
$ws.Range("H126").Value = 3230.8262  # was 3398.3157
$ws.Range("I126").Value = 3230.8262  # was 3398.3157
$ws.Range("K126").Value = 9692.4786  # was 10194.9471
$ws.Range("M126").Value = -7222.4786  # was -7724.947100000001

Write-Host "Applied 217 cell changes (214 updates, 2 clears, 1 new)."